# CORE_holdings.xlsx update: fill in the "Percent Change" column (E) for
# rows 2-8, which previously held placeholder zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected, so unlock it before writing, then restore
# protection afterwards (DrawingObjects/Contents/Scenarios protected,
# column & row formatting left allowed - matching the sheet's original
# protection intent).
$ws.Unprotect()

$ws.Range("E2").Value = 0.001628001628001696
$ws.Range("E3").Value = 0.005115778136779658
$ws.Range("E4").Value = 0.0004947433518862621
$ws.Range("E5").Value = 0.0011522949875169
$ws.Range("E6").Value = -0.008262052877138326
$ws.Range("E7").Value = -0.0044322176371695
$ws.Range("E8").Value = 0.002015579590251004

$ws.Protect($null, $true, $true, $true, $false, $false, $true, $true, $false, $false, $false, $false, $false, $false, $false, $false)
